$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new row of data (A8 = 0, B8 = "Unknown") right below the existing table rows
$ws.Range("A8").Value = 0
$ws.Range("B8").Value = "Unknown"

# Resize the table (ListObject) to include the new row
$table = $ws.ListObjects.Item("Table1")
$table.Resize($ws.Range("A1:B8"))

# Move the active selection to B9 (just below the newly extended table), matching the diff
$ws.Range("B9").Select()
